# ReportingOrganisationGroup: swap the "group-code" and "group-name" columns
# (column C <-> column D), including the header row, for every used row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $groupCode = $ws.Cells.Item($r, 3).Value2
    $groupName = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $groupName
    $ws.Cells.Item($r, 4).Value2 = $groupCode
}
